$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 639 (pushes old 639..653 down to 643..657),
# inheriting formatting (incl. the date number-format on column D) from row 638 above,
# matching native Excel "Insert" behaviour.
$ws.Rows("639:642").Insert()

# New row 639: Granny Smith / Primera
$ws.Range("A639").Value = 4
$ws.Range("B639").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C639").Value = "Los Lagos"
$ws.Range("D639").Value = "9/9/2021"
$ws.Range("E639").Value = 10
$ws.Range("F639").Value = "Fruta"
$ws.Range("G639").Value = 100104
$ws.Range("H639").Value = "Frutos de pepita"
$ws.Range("I639").Value = 100104002
$ws.Range("J639").Value = "Manzana"
$ws.Range("K639").Value = "Granny Smith"
$ws.Range("L639").Value = "Primera"
$ws.Range("M639").Value = 200
$ws.Range("N639").Value = 16000
$ws.Range("O639").Value = 16000
$ws.Range("P639").Value = 16000
$ws.Range("Q639").Value = "$/caja 16 kilos empedrada"
$ws.Range("R639").Value = "Provincia de Curicó"
$ws.Range("S639").Value = 1000
$ws.Range("T639").Value = 16

# New row 640: Granny Smith / Segunda
$ws.Range("A640").Value = 4
$ws.Range("B640").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C640").Value = "Los Lagos"
$ws.Range("D640").Value = "9/9/2021"
$ws.Range("E640").Value = 10
$ws.Range("F640").Value = "Fruta"
$ws.Range("G640").Value = 100104
$ws.Range("H640").Value = "Frutos de pepita"
$ws.Range("I640").Value = 100104002
$ws.Range("J640").Value = "Manzana"
$ws.Range("K640").Value = "Granny Smith"
$ws.Range("L640").Value = "Segunda"
$ws.Range("M640").Value = 100
$ws.Range("N640").Value = 11000
$ws.Range("O640").Value = 11000
$ws.Range("P640").Value = 11000
$ws.Range("Q640").Value = "$/caja 16 kilos empedrada"
$ws.Range("R640").Value = "Provincia de Curicó"
$ws.Range("S640").Value = 688
$ws.Range("T640").Value = 16

# New row 641: Royal Gala / Primera
$ws.Range("A641").Value = 4
$ws.Range("B641").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C641").Value = "Los Lagos"
$ws.Range("D641").Value = "9/9/2021"
$ws.Range("E641").Value = 10
$ws.Range("F641").Value = "Fruta"
$ws.Range("G641").Value = 100104
$ws.Range("H641").Value = "Frutos de pepita"
$ws.Range("I641").Value = 100104002
$ws.Range("J641").Value = "Manzana"
$ws.Range("K641").Value = "Royal Gala"
$ws.Range("L641").Value = "Primera"
$ws.Range("M641").Value = 200
$ws.Range("N641").Value = 16000
$ws.Range("O641").Value = 16000
$ws.Range("P641").Value = 16000
$ws.Range("Q641").Value = "$/caja 16 kilos empedrada"
$ws.Range("R641").Value = "Provincia de Curicó"
$ws.Range("S641").Value = 1000
$ws.Range("T641").Value = 16

# New row 642: Royal Gala / Segunda
$ws.Range("A642").Value = 4
$ws.Range("B642").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C642").Value = "Los Lagos"
$ws.Range("D642").Value = "9/9/2021"
$ws.Range("E642").Value = 10
$ws.Range("F642").Value = "Fruta"
$ws.Range("G642").Value = 100104
$ws.Range("H642").Value = "Frutos de pepita"
$ws.Range("I642").Value = 100104002
$ws.Range("J642").Value = "Manzana"
$ws.Range("K642").Value = "Royal Gala"
$ws.Range("L642").Value = "Segunda"
$ws.Range("M642").Value = 100
$ws.Range("N642").Value = 11000
$ws.Range("O642").Value = 11000
$ws.Range("P642").Value = 11000
$ws.Range("Q642").Value = "$/caja 16 kilos empedrada"
$ws.Range("R642").Value = "Provincia de Curicó"
$ws.Range("S642").Value = 688
$ws.Range("T642").Value = 16
